# Added Jeremy to people spreadsheet.
# Appends a new row (row 27) with Jeremy Dewar's contact info, mirroring
# the layout/formatting of the other "Post-Doctoral Researcher" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new row's values --------------------------------------
$ws.Range("B27").Value = "Jeremy Dewar"
$ws.Range("C27").Value = "jdewar@tulane.edu"
$ws.Range("E27").Value = "Post-Doctoral Researcher"
$ws.Range("F27").Value = "https://github.com/jdewar"
$ws.Range("G27").Value = "http://math.tulane.edu/~jdewar/"

# --- Match formatting used elsewhere on the sheet -----------------------
# Plain "filled" blank cells (same look as D4/G4/H4/I4/K4).
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").Copy() | Out-Null
$ws.Range("I27").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").Copy() | Out-Null
$ws.Range("J27").PasteSpecial(-4122) | Out-Null

# Plain text style used for the name/title/link columns (same as C26).
$ws.Range("C26").Copy() | Out-Null
$ws.Range("B27").PasteSpecial(-4122) | Out-Null
$ws.Range("C26").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("C26").Copy() | Out-Null
$ws.Range("F27").PasteSpecial(-4122) | Out-Null
$ws.Range("C26").Copy() | Out-Null
$ws.Range("G27").PasteSpecial(-4122) | Out-Null

# Email-style formatting for the email cell (same as J11).
$ws.Range("J11").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Match the row height used by similarly styled rows.
$ws.Rows(27).RowHeight = 15.75

# --- Restore the cursor/selection position left by the author -----------
$ws.Range("C20").Select() | Out-Null
